# ==========================================================================
# 20221111_B社予実データ.xlsx — "Add files via upload" commit replication
#
# Sheet1 ("B社　実績値input" = actual monthly sales figures) and
# Sheet2 ("B社　予測値input" = 12-month rolling sales forecast made each
# month) both get 12 new monthly rows appended (2021/04 .. 2022/03,
# rows 39-50). In addition, sheet2's most recent rows (35-38) had their
# "future" forecast columns still showing the "-" placeholder text — now
# that time has passed those forecasts are filled in with the actual
# numbers that were produced at the time.
# ==========================================================================

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # "B社　実績値input"
$ws2 = $wb.Worksheets.Item(2)   # "B社　予測値input"

# --- Fix sheet2 ("B社 予測値input") rows 35-38: replace placeholder cells with the now-known forecast values ---
$ws2.Range("J35").Value2 = 137000
$ws2.Range("K35").Value2 = 137000
$ws2.Range("L35").Value2 = 137000
$ws2.Range("M35").Value2 = 152000
$ws2.Range("N35").Value2 = 152000
$ws2.Range("I36").Value2 = 182000
$ws2.Range("J36").Value2 = 182000
$ws2.Range("K36").Value2 = 182000
$ws2.Range("L36").Value2 = 122000
$ws2.Range("M36").Value2 = 122000
$ws2.Range("N36").Value2 = 122000
$ws2.Range("H37").Value2 = 119000
$ws2.Range("I37").Value2 = 110000
$ws2.Range("J37").Value2 = 110000
$ws2.Range("K37").Value2 = 91000
$ws2.Range("L37").Value2 = 91000
$ws2.Range("M37").Value2 = 91000
$ws2.Range("N37").Value2 = 107000
$ws2.Range("G38").Value2 = 149000
$ws2.Range("H38").Value2 = 117000
$ws2.Range("I38").Value2 = 100000
$ws2.Range("J38").Value2 = 91000
$ws2.Range("K38").Value2 = 91000
$ws2.Range("L38").Value2 = 91000
$ws2.Range("M38").Value2 = 67000
$ws2.Range("N38").Value2 = 72000

# --- Sheet1 ("B社 実績値input"): append months 202104-202203 (rows 39-50) ---
$src1 = $ws1.Range("B38:C38")
$src1.Copy($ws1.Range("B39:C39"))
$ws1.Range("B39").Value2 = 202104
$ws1.Range("C39").Value2 = 171280
$src1.Copy($ws1.Range("B40:C40"))
$ws1.Range("B40").Value2 = 202105
$ws1.Range("C40").Value2 = 177500
$src1.Copy($ws1.Range("B41:C41"))
$ws1.Range("B41").Value2 = 202106
$ws1.Range("C41").Value2 = 199180
$src1.Copy($ws1.Range("B42:C42"))
$ws1.Range("B42").Value2 = 202107
$ws1.Range("C42").Value2 = 163180
$src1.Copy($ws1.Range("B43:C43"))
$ws1.Range("B43").Value2 = 202108
$ws1.Range("C43").Value2 = 134120
$src1.Copy($ws1.Range("B44:C44"))
$ws1.Range("B44").Value2 = 202109
$ws1.Range("C44").Value2 = 100920
$src1.Copy($ws1.Range("B45:C45"))
$ws1.Range("B45").Value2 = 202110
$ws1.Range("C45").Value2 = 48120
$src1.Copy($ws1.Range("B46:C46"))
$ws1.Range("B46").Value2 = 202111
$ws1.Range("C46").Value2 = 59080
$src1.Copy($ws1.Range("B47:C47"))
$ws1.Range("B47").Value2 = 202112
$ws1.Range("C47").Value2 = 44280
$src1.Copy($ws1.Range("B48:C48"))
$ws1.Range("B48").Value2 = 202201
$ws1.Range("C48").Value2 = 124980
$src1.Copy($ws1.Range("B49:C49"))
$ws1.Range("B49").Value2 = 202202
$ws1.Range("C49").Value2 = 118080
$src1.Copy($ws1.Range("B50:C50"))
$ws1.Range("B50").Value2 = 202203
$ws1.Range("C50").Value2 = 75000

# --- Sheet2 ("B社 予測値input"): append months 202104-202203 (rows 39-50) ---
$src2 = $ws2.Range("B38:N38")
$src2.Copy($ws2.Range("B39:N39"))
$ws2.Range("B39").Value2 = 202104
$ws2.Range("C39").Value2 = 175000
$ws2.Range("D39").Value2 = 205000
$ws2.Range("E39").Value2 = 185000
$ws2.Range("F39").Value2 = 148000
$ws2.Range("G39").Value2 = 115000
$ws2.Range("H39").Value2 = 100000
$ws2.Range("I39").Value2 = 59000
$ws2.Range("J39").Value2 = 60000
$ws2.Range("K39").Value2 = 59000
$ws2.Range("L39").Value2 = 101000
$ws2.Range("M39").Value2 = 100000
$ws2.Range("N39").Value2 = 101000
$src2.Copy($ws2.Range("B40:N40"))
$ws2.Range("B40").Value2 = 202105
$ws2.Range("C40").Value2 = 214200
$ws2.Range("D40").Value2 = 213200
$ws2.Range("E40").Value2 = 130200
$ws2.Range("F40").Value2 = 105800
$ws2.Range("G40").Value2 = 89800
$ws2.Range("H40").Value2 = 58800
$ws2.Range("I40").Value2 = 59800
$ws2.Range("J40").Value2 = 58800
$ws2.Range("K40").Value2 = 101000
$ws2.Range("L40").Value2 = 100000
$ws2.Range("M40").Value2 = 101000
$ws2.Range("N40").ClearContents()
$src2.Copy($ws2.Range("B41:N41"))
$ws2.Range("B41").Value2 = 202106
$ws2.Range("C41").Value2 = 227200
$ws2.Range("D41").Value2 = 180200
$ws2.Range("E41").Value2 = 88600
$ws2.Range("F41").Value2 = 79800
$ws2.Range("G41").Value2 = 58800
$ws2.Range("H41").Value2 = 59800
$ws2.Range("I41").Value2 = 58800
$ws2.Range("J41").Value2 = 101000
$ws2.Range("K41").Value2 = 100000
$ws2.Range("L41").Value2 = 101000
$ws2.Range("M41").ClearContents()
$ws2.Range("N41").ClearContents()
$src2.Copy($ws2.Range("B42:N42"))
$ws2.Range("B42").Value2 = 202107
$ws2.Range("C42").Value2 = 165600
$ws2.Range("D42").Value2 = 108400
$ws2.Range("E42").Value2 = 79400
$ws2.Range("F42").Value2 = 68800
$ws2.Range("G42").Value2 = 69800
$ws2.Range("H42").Value2 = 58800
$ws2.Range("I42").Value2 = 101000
$ws2.Range("J42").Value2 = 100000
$ws2.Range("K42").Value2 = 101000
$ws2.Range("L42").ClearContents()
$ws2.Range("M42").ClearContents()
$ws2.Range("N42").ClearContents()
$src2.Copy($ws2.Range("B43:N43"))
$ws2.Range("B43").Value2 = 202108
$ws2.Range("C43").Value2 = 146400
$ws2.Range("D43").Value2 = 69400
$ws2.Range("E43").Value2 = 58800
$ws2.Range("F43").Value2 = 59800
$ws2.Range("G43").Value2 = 53800
$ws2.Range("H43").Value2 = 121000
$ws2.Range("I43").Value2 = 120000
$ws2.Range("J43").Value2 = 136000
$ws2.Range("K43").ClearContents()
$ws2.Range("L43").ClearContents()
$ws2.Range("M43").ClearContents()
$ws2.Range("N43").ClearContents()
$src2.Copy($ws2.Range("B44:N44"))
$ws2.Range("B44").Value2 = 202109
$ws2.Range("C44").Value2 = 81400
$ws2.Range("D44").Value2 = 58800
$ws2.Range("E44").Value2 = 59800
$ws2.Range("F44").Value2 = 53800
$ws2.Range("G44").Value2 = 121000
$ws2.Range("H44").Value2 = 120000
$ws2.Range("I44").Value2 = 136000
$ws2.Range("J44").ClearContents()
$ws2.Range("K44").ClearContents()
$ws2.Range("L44").ClearContents()
$ws2.Range("M44").ClearContents()
$ws2.Range("N44").ClearContents()
$src2.Copy($ws2.Range("B45:N45"))
$ws2.Range("B45").Value2 = 202110
$ws2.Range("C45").Value2 = 38800
$ws2.Range("D45").Value2 = 59800
$ws2.Range("E45").Value2 = 53800
$ws2.Range("F45").Value2 = 121000
$ws2.Range("G45").Value2 = 120000
$ws2.Range("H45").Value2 = 136000
$ws2.Range("I45").ClearContents()
$ws2.Range("J45").ClearContents()
$ws2.Range("K45").ClearContents()
$ws2.Range("L45").ClearContents()
$ws2.Range("M45").ClearContents()
$ws2.Range("N45").ClearContents()
$src2.Copy($ws2.Range("B46:N46"))
$ws2.Range("B46").Value2 = 202111
$ws2.Range("C46").Value2 = 49800
$ws2.Range("D46").Value2 = 53800
$ws2.Range("E46").Value2 = 130000
$ws2.Range("F46").Value2 = 130000
$ws2.Range("G46").Value2 = 130000
$ws2.Range("H46").ClearContents()
$ws2.Range("I46").ClearContents()
$ws2.Range("J46").ClearContents()
$ws2.Range("K46").ClearContents()
$ws2.Range("L46").ClearContents()
$ws2.Range("M46").ClearContents()
$ws2.Range("N46").ClearContents()
$src2.Copy($ws2.Range("B47:N47"))
$ws2.Range("B47").Value2 = 202112
$ws2.Range("C47").Value2 = 44160
$ws2.Range("D47").Value2 = 129600
$ws2.Range("E47").Value2 = 104200
$ws2.Range("F47").Value2 = 94200
$ws2.Range("G47").Value2 = 133700
$ws2.Range("H47").Value2 = 133300
$ws2.Range("I47").Value2 = 133300
$ws2.Range("J47").Value2 = 133300
$ws2.Range("K47").Value2 = 133300
$ws2.Range("L47").Value2 = 133300
$ws2.Range("M47").Value2 = 145800
$ws2.Range("N47").Value2 = 120800
$src2.Copy($ws2.Range("B48:N48"))
$ws2.Range("B48").Value2 = 202201
$ws2.Range("C48").Value2 = 169600
$ws2.Range("D48").Value2 = 54200
$ws2.Range("E48").Value2 = 94200
$ws2.Range("F48").Value2 = 133700
$ws2.Range("G48").Value2 = 133300
$ws2.Range("H48").Value2 = 133300
$ws2.Range("I48").Value2 = 133300
$ws2.Range("J48").Value2 = 133300
$ws2.Range("K48").Value2 = 133300
$ws2.Range("L48").Value2 = 145800
$ws2.Range("M48").Value2 = 120800
$ws2.Range("N48").Value2 = 133300
$src2.Copy($ws2.Range("B49:N49"))
$ws2.Range("B49").Value2 = 202202
$ws2.Range("C49").Value2 = 118400
$ws2.Range("D49").Value2 = 75000
$ws2.Range("E49").Value2 = 133000
$ws2.Range("F49").Value2 = 133300
$ws2.Range("G49").Value2 = 133300
$ws2.Range("H49").Value2 = 133300
$ws2.Range("I49").Value2 = 133300
$ws2.Range("J49").Value2 = 133300
$ws2.Range("K49").Value2 = 133300
$ws2.Range("L49").Value2 = 133300
$ws2.Range("M49").Value2 = 133300
$ws2.Range("N49").Value2 = 133300
$src2.Copy($ws2.Range("B50:N50"))
$ws2.Range("B50").Value2 = 202203
$ws2.Range("C50").Value2 = 75000
$ws2.Range("D50").Value2 = 133000
$ws2.Range("E50").Value2 = 133300
$ws2.Range("F50").Value2 = 133300
$ws2.Range("G50").Value2 = 133300
$ws2.Range("H50").Value2 = 133300
$ws2.Range("I50").Value2 = 133300
$ws2.Range("J50").Value2 = 133300
$ws2.Range("K50").Value2 = 133300
$ws2.Range("L50").Value2 = 133300
$ws2.Range("M50").Value2 = 133300
$ws2.Range("N50").Value2 = 133300

# --- View state best-effort: sheet1 stays the active/selected sheet with
#     the cursor back at the top (A1); sheet2's cursor stays at B1 but the
#     window is scrolled down towards the newly-appended rows. ---
$ws2.Activate()
$excel.ActiveWindow.ScrollRow = 41
$ws2.Range("B1").Select()
$ws1.Activate()
$ws1.Range("A1").Select()
